# "Planilha investimento.xlsx" update:
#  - Salary bumped from 3500 to 4500 (and all dependent computed cells
#    recalculate: suggested investment, accumulated patrimony, dividends,
#    the 2/5/10/15/20/30-year projections, and the allocation table).
#  - Investment profile dropdown switched from MODERADO to CONSEVARDOR,
#    which re-drives the VLOOKUP-based allocation percentages/amounts.
#  - Label/value cells for the summary block are merged (B18:C18 ... B22:C22).
#  - The no-longer-needed chart helper defined names (_xlchart.v1.*) and the
#    stale external workbook link are removed.
#  - Minor view tweaks: headers hidden, view focused back on the
#    profile/allocation area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APP")

# --- Core input changes -----------------------------------------------
$ws.Range("D13").Value = 4500
$ws.Range("C32").Value = "CONSEVARDOR"

# --- Merge the label/value pairs that used to be single (unmerged) cells
$ws.Range("B18:C18").Merge()
$ws.Range("B19:C19").Merge()
$ws.Range("B20:C20").Merge()
$ws.Range("B21:C21").Merge()
$ws.Range("B22:C22").Merge()

# --- Remove the now-unused chart helper defined names ------------------
foreach ($n in @("_xlchart.v1.0", "_xlchart.v1.1", "_xlchart.v1.2")) {
    try {
        $wb.Names.Item($n).Delete()
    } catch {
    }
}

# --- Remove the stale external workbook reference ----------------------
$links = $wb.LinkSources()
if ($links) {
    foreach ($l in $links) {
        try {
            $wb.BreakLink($l, 1)
        } catch {
        }
    }
}

# --- View tweaks: hide row/col headers, return focus to the profile cell
$ws.Activate()
$excel.ActiveWindow.DisplayHeadings = $false
$ws.Range("C32").Select()

$wb.Application.CalculateFullRebuild()
